$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 74.40000000000001
$ws.Range("I15").Value = 74.40000000000001
$ws.Range("K15").Value = 223.2
$ws.Range("M15").Value = -54.20000000000002

$ws.Range("H29").Value = 955.6667
$ws.Range("J29").Value = 4000
$ws.Range("L29").Value = 12000
$ws.Range("N29").Value = -12562

$ws.Range("H40").Value = 999.6667
$ws.Range("I40").Value = 999.6667
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 999.6667
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -824.6667
$ws.Range("N40").ClearContents()

$ws.Range("H75").Value = 28400
$ws.Range("J75").Value = 28400
$ws.Range("L75").Value = 28400
$ws.Range("N75").Value = -30272

$ws.Range("H78").Value = 28400
$ws.Range("J78").Value = 28400
$ws.Range("L78").Value = 85200
$ws.Range("N78").Value = -94560

$ws.Range("H95").Value = 42500
$ws.Range("J95").Value = 42500
$ws.Range("L95").Value = 42500
$ws.Range("N95").Value = -47992

$ws.Range("H105").Value = 39800
$ws.Range("J105").Value = 39800
$ws.Range("L105").Value = 39800
$ws.Range("N105").Value = -46788

$ws.Range("H135").Value = 1517.25
$ws.Range("I135").Value = 1241.1428
$ws.Range("J135").Value = 3450
$ws.Range("K135").Value = 11170.2852
$ws.Range("L135").Value = 31050
$ws.Range("M135").Value = -8635.2852
$ws.Range("N135").Value = -36120

$ws.Range("H137").Value = 4708.9487
$ws.Range("J137").Value = 4968.7144
$ws.Range("L137").Value = 14906.1432
$ws.Range("N137").Value = -20006.1432

$ws.Range("H138").Value = 4053.3372
$ws.Range("I138").Value = 2398.4
$ws.Range("J138").Value = 4271.0923
$ws.Range("K138").Value = 7195.200000000001
$ws.Range("L138").Value = 12813.2769
$ws.Range("M138").Value = -2055.200000000001
$ws.Range("N138").Value = -23093.2769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13362.275
$ws.Range("I32").Value = 9517.293
$ws.Range("J32").Value = 18937.5
$ws.Range("K32").Value = 9517.293
$ws.Range("L32").Value = 18937.5
$ws.Range("M32").Value = -9230.293
$ws.Range("N32").Value = -19511.5

$ws.Range("H61").Value = 4752.75
$ws.Range("I61").Value = 4670.3335
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 4670.3335
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -4458.3335
$ws.Range("N61").Value = -5424

$ws.Range("H74").Value = 4755.857
$ws.Range("I74").Value = 5213.7
$ws.Range("J74").Value = 3611.25
$ws.Range("K74").Value = 5213.7
$ws.Range("L74").Value = 3611.25
$ws.Range("M74").Value = -4339.7
$ws.Range("N74").Value = -5359.25

$ws.Range("H77").Value = 4755.857
$ws.Range("I77").Value = 5213.7
$ws.Range("J77").Value = 3611.25
$ws.Range("K77").Value = 26068.5
$ws.Range("L77").Value = 18056.25
$ws.Range("M77").Value = -21700.5
$ws.Range("N77").Value = -26792.25

$ws.Range("H132").Value = 4947.6665
$ws.Range("I132").Value = 2066.6667
$ws.Range("J132").Value = 5908
$ws.Range("K132").Value = 6200.000100000001
$ws.Range("L132").Value = 17724
$ws.Range("M132").Value = -3670.000100000001
$ws.Range("N132").Value = -22784

$ws.Range("H136").Value = 4752.75
$ws.Range("I136").Value = 4670.3335
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 14011.0005
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -11461.0005
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 57700
$ws.Range("J56").Value = 57700
$ws.Range("L56").Value = 57700
$ws.Range("N56").Value = -59178

$ws.Range("H57").Value = 52500
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 100000
$ws.Range("K57").Value = 5000
$ws.Range("L57").Value = 100000
$ws.Range("M57").Value = -4280
$ws.Range("N57").Value = -101440

$ws.Range("H88").Value = 43950
$ws.Range("J88").Value = 43950
$ws.Range("L88").Value = 43950
$ws.Range("N88").Value = -44762

$ws.Range("H91").Value = 43950
$ws.Range("J91").Value = 43950
$ws.Range("L91").Value = 43950
$ws.Range("N91").Value = -46758

$ws.Range("H92").Value = 29250
$ws.Range("J92").Value = 29250
$ws.Range("L92").Value = 29250
$ws.Range("N92").Value = -34242

$ws.Range("H94").Value = 18520050
$ws.Range("I94").Value = 25001344
$ws.Range("K94").Value = 25001344
$ws.Range("M94").Value = -25000893

$ws.Range("H107").Value = 1608.3889
$ws.Range("I107").Value = 1528.1666
$ws.Range("K107").Value = 1528.1666
$ws.Range("M107").Value = 391.8334

$ws.Range("H132").Value = 54000
$ws.Range("J132").Value = 54000
$ws.Range("L132").Value = 54000
$ws.Range("N132").Value = -64120

$ws.Range("H134").Value = 3849.3142
$ws.Range("I134").Value = 2029.6666
$ws.Range("K134").Value = 6088.9998
$ws.Range("M134").Value = -3553.9998

$ws.Range("H135").Value = 48995
$ws.Range("J135").Value = 48995
$ws.Range("L135").Value = 48995
$ws.Range("N135").Value = -59135

$ws.Range("H136").Value = 52500
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 100000
$ws.Range("K136").Value = 5000
$ws.Range("L136").Value = 100000
$ws.Range("M136").Value = 100
$ws.Range("N136").Value = -110200

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 44005.453
$ws.Range("J138").Value = 44005.453
$ws.Range("L138").Value = 44005.453
$ws.Range("N138").Value = -54285.453

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 48607.6
$ws.Range("J140").Value = 48607.6
$ws.Range("L140").Value = 48607.6
$ws.Range("N140").Value = -58967.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4640.087
$ws.Range("I31").Value = 2079.3333
$ws.Range("J31").Value = 5543.8823
$ws.Range("K31").Value = 2079.3333
$ws.Range("L31").Value = 5543.8823
$ws.Range("M31").Value = -1784.3333
$ws.Range("N31").Value = -6133.8823

$ws.Range("H34").Value = 4640.087
$ws.Range("I34").Value = 2079.3333
$ws.Range("J34").Value = 5543.8823
$ws.Range("K34").Value = 2079.3333
$ws.Range("L34").Value = 5543.8823
$ws.Range("M34").Value = -1877.3333
$ws.Range("N34").Value = -5947.8823

$ws.Range("H95").Value = 39395
$ws.Range("J95").Value = 39395
$ws.Range("L95").Value = 39395
$ws.Range("N95").Value = -44887

$ws.Range("H106").Value = 34925
$ws.Range("J106").Value = 34925
$ws.Range("L106").Value = 34925
$ws.Range("N106").Value = -37449

$ws.Range("H134").Value = 7853.4443
$ws.Range("I134").Value = 7433.0586
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 22299.1758
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -19764.1758
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3622.5454
$ws.Range("J80").Value = 3826.1052
$ws.Range("L80").Value = 11478.3156
$ws.Range("N80").Value = -13350.3156

$ws.Range("H83").Value = 3622.5454
$ws.Range("J83").Value = 3826.1052
$ws.Range("L83").Value = 34434.9468
$ws.Range("N83").Value = -43794.9468

$ws.Range("H113").Value = 565.71185
$ws.Range("J113").Value = 527.1070999999999
$ws.Range("L113").Value = 1581.3213
$ws.Range("N113").Value = -5921.3213

$ws.Range("H122").Value = 2674.2952
$ws.Range("I122").Value = 843.6667
$ws.Range("J122").Value = 2991.1345
$ws.Range("K122").Value = 7593.0003
$ws.Range("L122").Value = 26920.2105
$ws.Range("M122").Value = -5143.0003
$ws.Range("N122").Value = -31820.2105

$ws.Range("H137").Value = 7743.077
$ws.Range("I137").Value = 3005.9412
$ws.Range("J137").Value = 16691
$ws.Range("K137").Value = 9017.8236
$ws.Range("L137").Value = 50073
$ws.Range("M137").Value = -3917.8236
$ws.Range("N137").Value = -60273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2611.7058
$ws.Range("I93").Value = 2455.4443
$ws.Range("J93").Value = 2787.5
$ws.Range("K93").Value = 2455.4443
$ws.Range("L93").Value = 2787.5
$ws.Range("M93").Value = -1207.4443
$ws.Range("N93").Value = -5283.5

$ws.Range("H136").Value = 4066.348
$ws.Range("I136").Value = 2109.4546
$ws.Range("J136").Value = 5860.1665
$ws.Range("K136").Value = 6328.3638
$ws.Range("L136").Value = 17580.4995
$ws.Range("M136").Value = -3778.3638
$ws.Range("N136").Value = -22680.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 11950
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 18900
$ws.Range("K39").Value = 5000
$ws.Range("L39").Value = 18900
$ws.Range("M39").Value = -4587
$ws.Range("N39").Value = -19726

$ws.Range("H81").Value = 20090280
$ws.Range("I81").Value = 21429580
$ws.Range("J81").Value = 780
$ws.Range("K81").Value = 42859160
$ws.Range("L81").Value = 1560
$ws.Range("M81").Value = -42858099
$ws.Range("N81").Value = -3682

$ws.Range("H84").Value = 20090280
$ws.Range("I84").Value = 21429580
$ws.Range("J84").Value = 780
$ws.Range("K84").Value = 214295800
$ws.Range("L84").Value = 7800
$ws.Range("M84").Value = -214290496
$ws.Range("N84").Value = -18408

$ws.Range("H103").Value = 35763.168
$ws.Range("J103").Value = 35763.168
$ws.Range("L103").Value = 35763.168
$ws.Range("N103").Value = -38107.168

$ws.Range("H136").Value = 6147
$ws.Range("I136").Value = 1922.4286
$ws.Range("K136").Value = 5767.2858
$ws.Range("M136").Value = -3217.2858
